$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextCell $ws 'D2' '47.055.29'
Set-TextCell $ws 'E2' '  +4.56%  '

# Row 3
Set-TextCell $ws 'D3' '2.500.37'
Set-TextCell $ws 'E3' '  +3.02%  '

# Row 5
Set-TextCell $ws 'D5' '323.16'
Set-TextCell $ws 'E5' '  +1.50%  '

# Row 6
Set-TextCell $ws 'D6' '104.96'
Set-TextCell $ws 'E6' '  +0.92%  '

# Row 7
Set-TextCell $ws 'E7' '  +0.86%  '

# Row 8
Set-TextCell $ws 'D8' '1.00'
Set-TextCell $ws 'E8' '  +0.03%  '

# Row 9
Set-TextCell $ws 'D9' '0.542'
Set-TextCell $ws 'E9' '  +2.10%  '

# Row 10
Set-TextCell $ws 'D10' '36.95'
Set-TextCell $ws 'E10' '  +3.47%  '

# Row 11
Set-TextCell $ws 'D11' '0.0812'
Set-TextCell $ws 'E11' '  +1.06%  '

# Row 12
Set-TextCell $ws 'E12' '  +0.53%  '

# Row 13
Set-TextCell $ws 'E13' '  -0.71%  '

# Row 14
Set-TextCell $ws 'D14' '7.22'
Set-TextCell $ws 'E14' '  +3.74%  '

# Row 15
Set-TextCell $ws 'D15' '2.891.12'
Set-TextCell $ws 'E15' '  +3.09%  '

# Row 16
Set-TextCell $ws 'D16' '2.513.82'
Set-TextCell $ws 'E16' '  +2.96%  '

# Row 17
Set-TextCell $ws 'E17' '  +1.00%  '

# Row 18
Set-TextCell $ws 'D18' '47.004.22'
Set-TextCell $ws 'E18' '  +4.74%  '

# Row 19
Set-TextCell $ws 'D19' '12.62'
Set-TextCell $ws 'E19' '  +2.15%  '

# Row 20
Set-TextCell $ws 'D20' '6.60'
Set-TextCell $ws 'E20' '  +4.03%  '

# Row 21
Set-TextCell $ws 'D21' '0.0₃0934'
Set-TextCell $ws 'E21' '  +1.49%  '

# Row 22
Set-TextCell $ws 'D22' '70.58'
Set-TextCell $ws 'E22' '  +2.43%  '

# Row 23
Set-TextCell $ws 'D23' '250.55'
Set-TextCell $ws 'E23' '  +2.85%  '

# Row 24
Set-TextCell $ws 'E24' '  +2.84%  '

# Row 25
Set-TextCell $ws 'D25' '2.55'
Set-TextCell $ws 'E25' '  +2.20%  '

# Row 26
Set-TextCell $ws 'D26' '26.20'
Set-TextCell $ws 'E26' '  +3.27%  '

# Row 27
Set-TextCell $ws 'D27' '1.00'
Set-TextCell $ws 'E27' '  -0.10%  '

# Row 28
Set-TextCell $ws 'E28' '  +4.69%  '

# Row 29
Set-TextCell $ws 'D29' '2.20'
Set-TextCell $ws 'E29' '  -2.95%  '

# Row 30
Set-TextCell $ws 'D30' '34.98'
Set-TextCell $ws 'E30' '  +3.01%  '

# Row 31
Set-TextCell $ws 'E31' '  +3.77%  '

# Row 32
Set-TextCell $ws 'D32' '49.48'
Set-TextCell $ws 'E32' '  +1.05%  '

# Row 33
Set-TextCell $ws 'D33' '19.61'
Set-TextCell $ws 'E33' '  -1.16%  '

# Row 34
Set-TextCell $ws 'E34' '  +1.55%  '

# Row 35
Set-TextCell $ws 'D35' '0.0777'
Set-TextCell $ws 'E35' '  +1.68%  '

# Row 36
Set-TextCell $ws 'E36' '  +0.20%  '

# Row 37
Set-TextCell $ws 'E37' '  +1.65%  '

# Row 38
Set-TextCell $ws 'D38' '4.57'
Set-TextCell $ws 'E38' '  +1.09%  '

# Row 39
Set-TextCell $ws 'D39' '2.96'
Set-TextCell $ws 'E39' '  +3.06%  '

# Row 40
Set-TextCell $ws 'B40' 'Monero'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D40' '122.49'
Set-TextCell $ws 'E40' '  -3.05%  '

# Row 41
Set-TextCell $ws 'B41' 'Stellar'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D41' '0.111'
Set-TextCell $ws 'E41' '  +1.34%  '

# Row 42
Set-TextCell $ws 'E42' '  +2.17%  '

# Row 43
Set-TextCell $ws 'D43' '21.28'
Set-TextCell $ws 'E43' '  +1.06%  '

# Row 44
Set-TextCell $ws 'D44' '0.0295'
Set-TextCell $ws 'E44' '  +1.88%  '

# Row 45
Set-TextCell $ws 'D45' '1.959.99'
Set-TextCell $ws 'E45' '  +0.77%  '

# Row 46
Set-TextCell $ws 'B46' 'ApeXProtocol'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws 'D46' '2.13'
Set-TextCell $ws 'E46' '  +0.71%  '

# Row 47
Set-TextCell $ws 'B47' 'NEARProtocol'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D47' '2.98'
Set-TextCell $ws 'E47' '  +0.72%  '

# Row 48
Set-TextCell $ws 'D48' '1.79'
Set-TextCell $ws 'E48' '  +0.44%  '

# Row 49
Set-TextCell $ws 'D49' '9.12'
Set-TextCell $ws 'E49' '  -1.18%  '

# Row 50
Set-TextCell $ws 'D50' '5.39'
Set-TextCell $ws 'E50' '  +15.10%  '

# Row 51
Set-TextCell $ws 'D51' '78.86'
Set-TextCell $ws 'E51' '  +4.22%  '
